$d = $word.ActiveDocument

# --- Title: change month from "أوت" (August) to "سبتمبر" (September) ---
$d.Content.Find.Execute("أوت", $true, $true, $false, $false, $false, $true, 1, $false, "سبتمبر", 2) | Out-Null

# --- Table 1: per-commune figures + totals ---
$tbl1 = $d.Tables.Item(1)
$tbl1.Cell(2, 5).Range.Text = "10 530 000,00"  # was 10 550 000,00
$tbl1.Cell(3, 3).Range.Text = "165"  # was 167
$tbl1.Cell(3, 4).Range.Text = "1 650 000,00"  # was 1 670 000,00
$tbl1.Cell(4, 3).Range.Text = "201"  # was 202
$tbl1.Cell(4, 4).Range.Text = "2 010 000,00"  # was 2 020 000,00
$tbl1.Cell(4, 5).Range.Text = "4 900 000,00"  # was 4 980 000,00
$tbl1.Cell(5, 3).Range.Text = "152"  # was 155
$tbl1.Cell(5, 4).Range.Text = "1 520 000,00"  # was 1 550 000,00
$tbl1.Cell(6, 3).Range.Text = "89"  # was 93
$tbl1.Cell(6, 4).Range.Text = "890 000,00"  # was 930 000,00
$tbl1.Cell(8, 3).Range.Text = "194"  # was 196
$tbl1.Cell(8, 4).Range.Text = "1 940 000,00"  # was 1 960 000,00
$tbl1.Cell(8, 5).Range.Text = "4 920 000,00"  # was 4 930 000,00
$tbl1.Cell(9, 3).Range.Text = "123"  # was 121
$tbl1.Cell(9, 4).Range.Text = "1 230 000,00"  # was 1 210 000,00
$tbl1.Cell(12, 3).Range.Text = "62"  # was 63
$tbl1.Cell(12, 4).Range.Text = "620 000,00"  # was 630 000,00
$tbl1.Cell(13, 3).Range.Text = "356"  # was 361
$tbl1.Cell(13, 4).Range.Text = "3 560 000,00"  # was 3 610 000,00
$tbl1.Cell(13, 5).Range.Text = "4 710 000,00"  # was 4 750 000,00
$tbl1.Cell(14, 3).Range.Text = "30"  # was 29
$tbl1.Cell(14, 4).Range.Text = "300 000,00"  # was 290 000,00
$tbl1.Cell(17, 3).Range.Text = "159"  # was 166
$tbl1.Cell(17, 4).Range.Text = "1 590 000,00"  # was 1 660 000,00
$tbl1.Cell(17, 5).Range.Text = "3 900 000,00"  # was 3 980 000,00
$tbl1.Cell(18, 3).Range.Text = "26"  # was 27
$tbl1.Cell(18, 4).Range.Text = "260 000,00"  # was 270 000,00
$tbl1.Cell(19, 3).Range.Text = "105"  # was 103
$tbl1.Cell(19, 4).Range.Text = "1 050 000,00"  # was 1 030 000,00
$tbl1.Cell(20, 3).Range.Text = "100"  # was 102
$tbl1.Cell(20, 4).Range.Text = "1 000 000,00"  # was 1 020 000,00
$tbl1.Cell(21, 5).Range.Text = "3 750 000,00"  # was 3 770 000,00
$tbl1.Cell(22, 3).Range.Text = "130"  # was 134
$tbl1.Cell(22, 4).Range.Text = "1 300 000,00"  # was 1 340 000,00
$tbl1.Cell(23, 3).Range.Text = "96"  # was 94
$tbl1.Cell(23, 4).Range.Text = "960 000,00"  # was 940 000,00
$tbl1.Cell(25, 3).Range.Text = "3271"  # was 3296
$tbl1.Cell(25, 4).Range.Text = "32 710 000,00"  # was 32 960 000,00
$tbl1.Cell(25, 5).Range.Text = "32 710 000,00"  # was 32 960 000,00

# --- Table 2: per-commune figures + totals (second district block) ---
$tbl2 = $d.Tables.Item(2)
$tbl2.Cell(2, 3).Range.Text = "477"  # was 480
$tbl2.Cell(2, 4).Range.Text = "4 770 000,00"  # was 4 800 000,00
$tbl2.Cell(2, 5).Range.Text = "6 140 000,00"  # was 6 180 000,00
$tbl2.Cell(3, 3).Range.Text = "83"  # was 84
$tbl2.Cell(3, 4).Range.Text = "830 000,00"  # was 840 000,00
$tbl2.Cell(5, 3).Range.Text = "204"  # was 209
$tbl2.Cell(5, 4).Range.Text = "2 040 000,00"  # was 2 090 000,00
$tbl2.Cell(5, 5).Range.Text = "2 710 000,00"  # was 2 800 000,00
$tbl2.Cell(6, 3).Range.Text = "67"  # was 71
$tbl2.Cell(6, 4).Range.Text = "670 000,00"  # was 710 000,00
$tbl2.Cell(7, 3).Range.Text = "885"  # was 898
$tbl2.Cell(7, 4).Range.Text = "8 850 000,00"  # was 8 980 000,00
$tbl2.Cell(7, 5).Range.Text = "8 850 000,00"  # was 8 980 000,00
$tbl2.Cell(8, 3).Range.Text = "4156"  # was 4194
$tbl2.Cell(8, 4).Range.Text = "41 560 000,00"  # was 41 940 000,00
$tbl2.Cell(8, 5).Range.Text = "41 560 000,00"  # was 41 940 000,00

# --- Spelled-out grand total in Arabic words ---
$d.Content.Find.Execute("واحد وأربعون مليون وتسعمئة وأربعون ألف", $true, $true, $false, $false, $false, $true, 1, $false, "واحد وأربعون مليون وخمسمئة وستون ألف", 2) | Out-Null
